$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the literal text into the cell even when it looks like a number
    # (e.g. "0.9994", "241.50") so the cell keeps matching the source
    # sheet's text representation instead of becoming a numeric cell.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '29.957.67'
$ws.Range('E2').Value = '  +0.70%  '
Set-TextValue $ws.Range('D3') '1.907.69'
$ws.Range('E3').Value = '  +1.03%  '
Set-TextValue $ws.Range('D4') '0.9994'
Set-TextValue $ws.Range('D5') '0.8102'
$ws.Range('E5').Value = '  +7.86%  '
Set-TextValue $ws.Range('D6') '241.50'
$ws.Range('E6').Value = '  +1.03%  '
Set-TextValue $ws.Range('D7') '0.9994'
$ws.Range('E7').Value = '  -0.17%  '
Set-TextValue $ws.Range('D8') '0.3125'
$ws.Range('E8').Value = '  +3.17%  '
Set-TextValue $ws.Range('D9') '26.41'
$ws.Range('E9').Value = '  +4.86%  '
Set-TextValue $ws.Range('D10') '0.06997'
$ws.Range('E10').Value = '  +3.02%  '
Set-TextValue $ws.Range('D11') '0.08007'
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D12') '0.7443'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D13') '1.903.89'
$ws.Range('E13').Value = '  +0.92%  '
Set-TextValue $ws.Range('D14') '5.189'
$ws.Range('E14').Value = '  +1.03%  '
Set-TextValue $ws.Range('D15') '92.54'
$ws.Range('E15').Value = '  +2.55%  '
Set-TextValue $ws.Range('D16') '29.960.79'
$ws.Range('E16').Value = '  +0.68%  '
Set-TextValue $ws.Range('D17') '14.00'
$ws.Range('E17').Value = '  +1.22%  '
Set-TextValue $ws.Range('D18') '5.872'
$ws.Range('E18').Value = '  -0.47%  '
Set-TextValue $ws.Range('D19') '245.42'
$ws.Range('E19').Value = '  +1.64%  '
Set-TextValue $ws.Range('D20') '0.000007785'
$ws.Range('E20').Value = '  +1.98%  '
Set-TextValue $ws.Range('D21') '1.001'
$ws.Range('E21').Value = '  -0.02%  '
Set-TextValue $ws.Range('D22') '2.152.24'
$ws.Range('E22').Value = '  +1.10%  '
Set-TextValue $ws.Range('D23') '0.9998'
$ws.Range('E23').Value = '  -0.16%  '
Set-TextValue $ws.Range('D24') '6.943'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D25') '0.1540'
$ws.Range('E25').Value = '  +21.13%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D26') '168.73'
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D27') '9.217'
$ws.Range('E27').Value = '  +0.53%  '
Set-TextValue $ws.Range('D28') '18.86'
$ws.Range('E28').Value = '  +1.46%  '
$ws.Range('E29').Value = '  +2.85%  '
$ws.Range('E30').Value = '  -1.62%  '
Set-TextValue $ws.Range('D31') '1.512'
$ws.Range('E31').Value = '  +0.01%  '
Set-TextValue $ws.Range('D32') '4.298'
$ws.Range('E32').Value = '  +1.66%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D33') '4.070'
$ws.Range('E33').Value = '  +1.88%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D34') '0.05531'
$ws.Range('E34').Value = '  +6.27%  '
Set-TextValue $ws.Range('D35') '1.263'
$ws.Range('E35').Value = '  +1.48%  '
Set-TextValue $ws.Range('D36') '0.7308'
$ws.Range('E36').Value = '  +1.06%  '
Set-TextValue $ws.Range('D37') '2.711'
$ws.Range('E37').Value = '  +0.11%  '
Set-TextValue $ws.Range('D38') '0.01919'
$ws.Range('E38').Value = '  +1.07%  '
Set-TextValue $ws.Range('D39') '2.790'
$ws.Range('E39').Value = '  +1.07%  '
Set-TextValue $ws.Range('D40') '0.4407'
$ws.Range('E40').Value = '  +0.93%  '
Set-TextValue $ws.Range('D41') '72.12'
$ws.Range('E41').Value = '  +1.57%  '
Set-TextValue $ws.Range('D42') '5.992'
$ws.Range('E42').Value = '  -2.10%  '
Set-TextValue $ws.Range('D43') '0.9996'
$ws.Range('E43').Value = '  -0.13%  '
Set-TextValue $ws.Range('D44') '0.8385'
$ws.Range('E44').Value = '  +1.60%  '
Set-TextValue $ws.Range('D45') '1.890'
$ws.Range('E45').Value = '  +0.91%  '
Set-TextValue $ws.Range('D46') '101.03'
$ws.Range('E46').Value = '  +1.77%  '
Set-TextValue $ws.Range('D47') '7.572'
$ws.Range('E47').Value = '  +0.27%  '
Set-TextValue $ws.Range('D48') '9.718'
$ws.Range('E48').Value = '  +0.66%  '
Set-TextValue $ws.Range('D49') '982.71'
$ws.Range('E49').Value = '  +9.85%  '
Set-TextValue $ws.Range('D50') '2.057.76'
$ws.Range('E50').Value = '  +0.58%  '
Set-TextValue $ws.Range('D51') '36.19'
$ws.Range('E51').Value = '  +1.02%  '
